$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B42").Value = 44022
$ws.Range("B42").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C42").Value = 35419
$ws.Range("D42").Value = 1230
$ws.Range("E42").Value = 16811
$ws.Range("F42").Value = 616
$ws.Range("G42").Value = 47.46
$ws.Range("H42").Value = 50.08

$ws.Range("I42").Value = $true
$ws.Range("J42").Value = $true

$ws.Range("O42").Value = "Success!"
